$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Clcf1"
$ws.Range("C2").Value = "Crlf1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.758258666666666
$ws.Range("H2").Value = 5.274775999999999
$ws.Range("I2").Value = 0.1132279568112417
$ws.Range("J2").Value = 0.1132279568112417
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1153933333333333
$ws.Range("N2").Value = 0.34618
$ws.Range("O2").Value = 0.007294522992813531
$ws.Range("P2").Value = 0.007294522992813531
$ws.Range("Q2").Value = 0.2028913284088888
$ws.Range("R2").Value = 1.82602195568
$ws.Range("S2").Value = 0.0008259439343888998
$ws.Range("T2").Value = 0.0008259439343888999

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Clcf1"
$ws.Range("C3").Value = "Crlf1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.758258666666666
$ws.Range("H3").Value = 5.274775999999999
$ws.Range("I3").Value = 0.1132279568112417
$ws.Range("J3").Value = 0.1132279568112417
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.096090333333334
$ws.Range("N3").Value = 24.288271
$ws.Range("O3").Value = 0.5117896795458609
$ws.Range("P3").Value = 0.5117896795458609
$ws.Range("Q3").Value = 14.23502099469956
$ws.Range("R3").Value = 128.115188952296
$ws.Range("S3").Value = 0.05794889973205795
$ws.Range("T3").Value = 0.05794889973205795

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Clcf1"
$ws.Range("C4").Value = "Crlf1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.758258666666666
$ws.Range("H4").Value = 5.274775999999999
$ws.Range("I4").Value = 0.1132279568112417
$ws.Range("J4").Value = 0.1132279568112417
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.447211
$ws.Range("N4").Value = 22.341633
$ws.Range("O4").Value = 0.4707711468470205
$ws.Range("P4").Value = 0.4707711468470205
$ws.Range("Q4").Value = 13.09412328324533
$ws.Range("R4").Value = 117.847109549208
$ws.Range("S4").Value = 0.05330445508317315
$ws.Range("T4").Value = 0.05330445508317316

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Clcf1"
$ws.Range("C5").Value = "Crlf1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.758258666666666
$ws.Range("H5").Value = 5.274775999999999
$ws.Range("I5").Value = 0.1132279568112417
$ws.Range("J5").Value = 0.1132279568112417
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.16048
$ws.Range("N5").Value = 0.48144
$ws.Range("O5").Value = 0.01014465061430512
$ws.Range("P5").Value = 0.01014465061430512
$ws.Range("Q5").Value = 0.2821653508266666
$ws.Range("R5").Value = 2.53948815744
$ws.Range("S5").Value = 0.001148658061621676
$ws.Range("T5").Value = 0.001148658061621677

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Clcf1"
$ws.Range("C6").Value = "Crlf1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.690388333333333
$ws.Range("H6").Value = 8.071165
$ws.Range("I6").Value = 0.1732550390834427
$ws.Range("J6").Value = 0.1732550390834427
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1153933333333333
$ws.Range("N6").Value = 0.34618
$ws.Range("O6").Value = 0.007294522992813531
$ws.Range("P6").Value = 0.007294522992813531
$ws.Range("Q6").Value = 0.3104528777444445
$ws.Range("R6").Value = 2.7940758997
$ws.Range("S6").Value = 0.00126381286621498
$ws.Range("T6").Value = 0.00126381286621498

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Clcf1"
$ws.Range("C7").Value = "Crlf1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.690388333333333
$ws.Range("H7").Value = 8.071165
$ws.Range("I7").Value = 0.1732550390834427
$ws.Range("J7").Value = 0.1732550390834427
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.096090333333334
$ws.Range("N7").Value = 24.288271
$ws.Range("O7").Value = 0.5117896795458609
$ws.Range("P7").Value = 0.5117896795458609
$ws.Range("Q7").Value = 21.78162697841278
$ws.Range("R7").Value = 196.034642805715
$ws.Range("S7").Value = 0.08867014093222074
$ws.Range("T7").Value = 0.08867014093222074

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Clcf1"
$ws.Range("C8").Value = "Crlf1"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.690388333333333
$ws.Range("H8").Value = 8.071165
$ws.Range("I8").Value = 0.1732550390834427
$ws.Range("J8").Value = 0.1732550390834427
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.447211
$ws.Range("N8").Value = 22.341633
$ws.Range("O8").Value = 0.4707711468470205
$ws.Range("P8").Value = 0.4707711468470205
$ws.Range("Q8").Value = 20.03588959027167
$ws.Range("R8").Value = 180.323006312445
$ws.Range("S8").Value = 0.08156347344633767
$ws.Range("T8").Value = 0.08156347344633769

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Clcf1"
$ws.Range("C9").Value = "Crlf1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.690388333333333
$ws.Range("H9").Value = 8.071165
$ws.Range("I9").Value = 0.1732550390834427
$ws.Range("J9").Value = 0.1732550390834427
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.16048
$ws.Range("N9").Value = 0.48144
$ws.Range("O9").Value = 0.01014465061430512
$ws.Range("P9").Value = 0.01014465061430512
$ws.Range("Q9").Value = 0.4317535197333333
$ws.Range("R9").Value = 3.8857816776
$ws.Range("S9").Value = 0.001757611838669304
$ws.Range("T9").Value = 0.001757611838669305

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Clcf1"
$ws.Range("C10").Value = "Crlf1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.938311666666667
$ws.Range("H10").Value = 29.814935
$ws.Range("I10").Value = 0.6400052196548212
$ws.Range("J10").Value = 0.640005219654821
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1153933333333333
$ws.Range("N10").Value = 0.34618
$ws.Range("O10").Value = 0.007294522992813531
$ws.Range("P10").Value = 0.007294522992813531
$ws.Range("Q10").Value = 1.146814910922222
$ws.Range("R10").Value = 10.3213341983
$ws.Range("S10").Value = 0.004668532790292767
$ws.Range("T10").Value = 0.004668532790292767

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Clcf1"
$ws.Range("C11").Value = "Crlf1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9.938311666666667
$ws.Range("H11").Value = 29.814935
$ws.Range("I11").Value = 0.6400052196548212
$ws.Range("J11").Value = 0.640005219654821
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.096090333333334
$ws.Range("N11").Value = 24.288271
$ws.Range("O11").Value = 0.5117896795458609
$ws.Range("P11").Value = 0.5117896795458609
$ws.Range("Q11").Value = 80.4614690141539
$ws.Range("R11").Value = 724.1532211273851
$ws.Range("S11").Value = 0.3275480662748192
$ws.Range("T11").Value = 0.3275480662748191

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Clcf1"
$ws.Range("C12").Value = "Crlf1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9.938311666666667
$ws.Range("H12").Value = 29.814935
$ws.Range("I12").Value = 0.6400052196548212
$ws.Range("J12").Value = 0.640005219654821
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 7.447211
$ws.Range("N12").Value = 22.341633
$ws.Range("O12").Value = 0.4707711468470205
$ws.Range("P12").Value = 0.4707711468470205
$ws.Range("Q12").Value = 74.01270396542834
$ws.Range("R12").Value = 666.1143356888551
$ws.Range("S12").Value = 0.3012959912449794
$ws.Range("T12").Value = 0.3012959912449794

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Clcf1"
$ws.Range("C13").Value = "Crlf1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 9.938311666666667
$ws.Range("H13").Value = 29.814935
$ws.Range("I13").Value = 0.6400052196548212
$ws.Range("J13").Value = 0.640005219654821
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.16048
$ws.Range("N13").Value = 0.48144
$ws.Range("O13").Value = 0.01014465061430512
$ws.Range("P13").Value = 0.01014465061430512
$ws.Range("Q13").Value = 1.594900256266667
$ws.Range("R13").Value = 14.3541023064
$ws.Range("S13").Value = 0.006492629344729764
$ws.Range("T13").Value = 0.006492629344729764

$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Clcf1"
$ws.Range("C14").Value = "Crlf1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.141526666666667
$ws.Range("H14").Value = 3.42458
$ws.Range("I14").Value = 0.0735117844504946
$ws.Range("J14").Value = 0.07351178445049458
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1153933333333333
$ws.Range("N14").Value = 0.34618
$ws.Range("O14").Value = 0.007294522992813531
$ws.Range("P14").Value = 0.007294522992813531
$ws.Range("Q14").Value = 0.1317245671555556
$ws.Range("R14").Value = 1.1855211044
$ws.Range("S14").Value = 0.000536233401916885
$ws.Range("T14").Value = 0.000536233401916885

$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Clcf1"
$ws.Range("C15").Value = "Crlf1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.141526666666667
$ws.Range("H15").Value = 3.42458
$ws.Range("I15").Value = 0.0735117844504946
$ws.Range("J15").Value = 0.07351178445049458
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 8.096090333333334
$ws.Range("N15").Value = 24.288271
$ws.Range("O15").Value = 0.5117896795458609
$ws.Range("P15").Value = 0.5117896795458609
$ws.Range("Q15").Value = 9.241903011242222
$ws.Range("R15").Value = 83.17712710118
$ws.Range("S15").Value = 0.03762257260676303
$ws.Range("T15").Value = 0.03762257260676302

$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Clcf1"
$ws.Range("C16").Value = "Crlf1"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.141526666666667
$ws.Range("H16").Value = 3.42458
$ws.Range("I16").Value = 0.0735117844504946
$ws.Range("J16").Value = 0.07351178445049458
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 7.447211
$ws.Range("N16").Value = 22.341633
$ws.Range("O16").Value = 0.4707711468470205
$ws.Range("P16").Value = 0.4707711468470205
$ws.Range("Q16").Value = 8.501189948793334
$ws.Range("R16").Value = 76.51070953914
$ws.Range("S16").Value = 0.03460722707253031
$ws.Range("T16").Value = 0.03460722707253031

$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Clcf1"
$ws.Range("C17").Value = "Crlf1"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.141526666666667
$ws.Range("H17").Value = 3.42458
$ws.Range("I17").Value = 0.0735117844504946
$ws.Range("J17").Value = 0.07351178445049458
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.16048
$ws.Range("N17").Value = 0.48144
$ws.Range("O17").Value = 0.01014465061430512
$ws.Range("P17").Value = 0.01014465061430512
$ws.Range("Q17").Value = 0.1831921994666666
$ws.Range("R17").Value = 1.6487297952
$ws.Range("S17").Value = 0.0007457513692843754
$ws.Range("T17").Value = 0.0007457513692843754

